$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 45/46 swap places (Stellar now ranks above ApeXProtocol); update the
# coin name and coinranking.com link columns to match the new ordering.
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"

# Column D (Price) holds numeric-looking values stored as plain text. Force the
# cells to Text format before writing so Excel does not silently convert them
# to floating point numbers (which would corrupt values like "53.70" or "0.140").
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.816.08"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "3.501.90"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "597.43"
$ws.Range("E5").Value = "  -2.00%  "

$ws.Range("D6").Value = "194.37"
$ws.Range("E6").Value = "  +4.60%  "

$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -2.25%  "

$ws.Range("D10").Value = "0.649"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").Value = "53.70"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "0.0000301"
$ws.Range("E12").Value = "  -2.71%  "

$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "4.053.59"
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").Value = "607.05"
$ws.Range("E15").Value = "  +3.92%  "

$ws.Range("D16").Value = "69.947.57"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "18.98"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").Value = "3.488.72"
$ws.Range("E19").Value = "  -1.33%  "

$ws.Range("D20").Value = "0.121"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D22").Value = "17.89"
$ws.Range("E22").Value = "  +2.16%  "

$ws.Range("D23").Value = "104.16"
$ws.Range("E23").Value = "  +8.54%  "

$ws.Range("D24").Value = "5.14"
$ws.Range("E24").Value = "  +5.51%  "

$ws.Range("E25").Value = "  -3.05%  "

$ws.Range("D26").Value = "3.05"
$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +2.02%  "

$ws.Range("D29").Value = "33.53"
$ws.Range("E29").Value = "  +4.47%  "

$ws.Range("D30").Value = "4.57"
$ws.Range("E30").Value = "  +26.32%  "

$ws.Range("D31").Value = "7.06"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "12.61"
$ws.Range("E32").Value = "  +3.57%  "

$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").Value = "63.24"
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").Value = "3.732.39"
$ws.Range("E35").Value = "  +5.81%  "

$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("E36").Value = "  +4.69%  "

$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  -6.70%  "

$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("E40").Value = "  +1.69%  "

$ws.Range("D41").Value = "36.59"
$ws.Range("E41").Value = "  -1.54%  "

$ws.Range("D42").Value = "498.32"
$ws.Range("E42").Value = "  -7.51%  "

$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").Value = "0.0456"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").Value = "0.140"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("D46").Value = "3.32"
$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").Value = "8.71"
$ws.Range("E49").Value = "  -4.66%  "

$ws.Range("D50").Value = "131.78"
$ws.Range("E50").Value = "  -2.10%  "

$ws.Range("E51").Value = "  +0.01%  "

# Restore the default (unformatted/"General") style on column D now that the
# text values are set, matching the workbook's original formatting.
$priceRange.Style = "Normal"
